$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New iteration-4 block (rows 63-72) -------------------------------
# Row 63: date / iteration-number / header text, styled like the
# existing iteration-3 header row (54).
$ws.Cells.Item(54, 1).Copy()
$ws.Cells.Item(63, 1).PasteSpecial(-4122)
$ws.Cells.Item(63, 1).Value = 45663
$ws.Cells.Item(63, 2).Value = 4
$ws.Cells.Item(63, 3).Value = "Si identificano 6 casi d'uso da dettagliare"

# Row 72 typed next (matches the original authoring order captured by
# the shared-strings table order in the target workbook).
$ws.Cells.Item(72, 3).Value = "Associa Poi"

# Row 64
$ws.Cells.Item(64, 3).Value = "Inserisci nuova attività"
$ws.Cells.Item(64, 4).Value = 1
$ws.Cells.Item(64, 6).Value = "dettagliato con diagramma"

# Row 66
$ws.Cells.Item(66, 3).Value = "Modifica Poi"
$ws.Cells.Item(66, 4).Value = 3
$ws.Cells.Item(66, 6).Value = "dettagliato con diagramma"

# Row 67
$ws.Cells.Item(67, 3).Value = "Rimuovi Poi"
$ws.Cells.Item(67, 4).Value = 4
$ws.Cells.Item(67, 6).Value = "dettagliato con diagramma"

# Row 68
$ws.Cells.Item(68, 3).Value = "Modifica Attività "
$ws.Cells.Item(68, 4).Value = 5

# Row 69
$ws.Cells.Item(69, 3).Value = "Rimuovi Attività"
$ws.Cells.Item(69, 4).Value = 6

# --- Rename the old "1 caso d'uso inserisci comune" use case -----------
$ws.Cells.Item(29, 3).Value = "Crea Nuovo Comune"

# Row 65 filled in last.
$ws.Cells.Item(65, 3).Value = "Ricerca Poi"
$ws.Cells.Item(65, 4).Value = 2
$ws.Cells.Item(65, 6).Value = "dettagliato con diagramma"

# --- Page setup: portrait / A4-ish "9" paper, like the target ----------
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1

# --- Update the on-screen selection to match the author's final spot ---
$ws.Range("F67").Select()
